# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" column (E) on the account-statement table listed the
# periods in descending order (2106, 2105, 2104, 2103). The refreshed data
# extract lists them in ascending order (2103, 2104, 2105, 2106), so update
# each row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("E16").Value = "2103"
$ws.Range("E17").Value = "2104"
$ws.Range("E18").Value = "2105"
$ws.Range("E19").Value = "2106"
